# Update the Task table: re-sort existing rows (IDs 1-3 reordered, Project D's
# task renamed/re-valued) and append four new Project D task rows, then grow
# the table/autofilter to match and refresh the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ID 1 / Project A / Design UI) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Project A"
$ws.Range("D2").Value = "Design UI"
$ws.Range("E2").Value = 100
$ws.Range("F2").ClearContents()
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 1234

# --- Row 3 (ID 2 / Project A / Implement Backend) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Project A"
$ws.Range("D3").Value = "Implement Backend"
$ws.Range("E3").Value = 55
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 1235

# --- Row 4 (ID 3 / Project A / Testing) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Project A"
$ws.Range("D4").Value = "Testing"
$ws.Range("E4").Value = 15
$ws.Range("I4").Value = 123

# --- Row 5 (ID 4 / Project B / API Development) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Project B"
$ws.Range("D5").Value = "API Development"
$ws.Range("E5").Value = 50
$ws.Range("F5").ClearContents()
$ws.Range("I5").Value = 1236

# Rows 6-8 are unchanged:
#  6: ID 5 / Project B / Database Setup / 60 / 1327
#  7: ID 6 / Project C / Database Setup / 30 / 35413
#  8: ID 7 / Project C / UI Development / 56 / 4654

# --- Row 9 (ID 8 / Project D / Analysis - renamed from Implement Backend) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Project D"
$ws.Range("D9").Value = "Analysis"
$ws.Range("E9").Value = 45
$ws.Range("I9").Value = 3135

# --- New row 10 (ID 9 / Project D / Database Setup) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Project D"
$ws.Range("D10").Value = "Database Setup"
$ws.Range("E10").Value = 75
$ws.Range("I10").Value = 1314

# --- New row 11 (ID 10 / Project D / Implement Backend) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "Project D"
$ws.Range("D11").Value = "Implement Backend"
$ws.Range("E11").Value = 47
$ws.Range("I11").Value = 13486

# --- New row 12 (ID 11 / Project D / UI Development) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = "Project D"
$ws.Range("D12").Value = "UI Development"
$ws.Range("E12").Value = 46
$ws.Range("I12").Value = 1314

# --- New row 13 (ID 12 / Project D / Testing) ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "Project D"
$ws.Range("D13").Value = "Testing"
$ws.Range("E13").Value = 47
$ws.Range("I13").Value = 13486

# Grow the table (and its autofilter) to cover the newly added rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I13"))

# Refresh the view: zoom in a bit and select the whole table like the
# recorded session did.
$excel.ActiveWindow.Zoom = 85
$null = $ws.Range("A1:I13").Select()
